$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.249.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.99%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.517.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.86%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.65%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.59%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.14"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.079.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.519.34"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.44"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.85%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.259.25"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.02"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.08%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +9.54%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.40"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.86%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.56%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.06"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.47"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "604.84"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.20%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.80"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.21%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.72"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.155"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0805"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.66%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.18%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.02"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.21%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.387"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.46%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.246.62"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.33%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0423"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.90%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.68"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.66"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.54"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.71%  "

# Row 45 (was ApeXProtocol -> now Fetch.AI)
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.32%  "

# Row 46 (was Fetch.AI -> now ApeXProtocol)
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.27%  "
